# Update the folder/file-path text in the "inputs_outputs" sheet so the
# tool points at the new (W:\ network-drive) locations instead of the old
# local (C:\Users\dpere\...) paths, and refresh the project sub-folder name
# plus the shapefile's embedded date stamp. Also move the active selection
# the way the author left it when saving (C3 instead of B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "location of base-forecast software" -> new network path
$ws.Range("B2").Value = "W:\Data\Forecast\Tools\forecast_git\create_forecast_basic\current"

# "output location by version" -> new network path
$ws.Range("B3").Value = "W:\Data\Forecast\forecast_by_version\V4\BASE_YEAR"

# "output location" -> new project folder
$ws.Range("B4").Value = "W:\Projects\תכניות מרחביות\דרום_מערב_122\קבצי עבודה\תחזיות_דמוגרפיות\For_approval\Reference_tabels"

# "new layer location" -> new project folder + refreshed shapefile date stamp
$ws.Range("B6").Value = "W:\Projects\תכניות מרחביות\דרום_מערב_122\קבצי עבודה\תחזיות_דמוגרפיות\For_approval\Reference_tabels\shp\TAZ_V4_241216_with_geo_info.shp"

# Leave the selection where the author left it on save.
$ws.Range("C3").Select()
